# The deck's "datetimeFigureOut" date placeholder (master + all 11 layouts)
# was re-cached from 9/19/2021 to 9/25/2021 (the file was re-saved six days
# later). Update the cached text of every "Date Placeholder" shape on the
# slide master and on each slide layout to match.

$p = $ppt.ActivePresentation
$newDate = "9/25/2021"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -eq "9/19/2021") {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout attached to the master
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}
